# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund-holding detail, same 8-column layout
# as the existing "2021-Q3" sheet) positioned right before the "总计"
# (totals) sheet, and prepends a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by duplicating "总计" (so it inherits
#    the same header/index-column styling) and inserting it immediately
#    before "总计" in the tab order.
# ---------------------------------------------------------------------
$totalSheet.Copy($totalSheet, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# Extend the header formatting (style copied along with "总计") from D1
# out to E1:H1 so all 8 header cells share the same look.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

# Add row 3 (second data row) by copying row 2's index-column (A) format.
$newSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row 2 — leading apostrophe keeps numeric-looking text (fund code /
# percentages / amounts) stored as text instead of being coerced to a number.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'009726"
$newSheet.Range("C2").Value = "招商中证500等权重指数增强A"
$newSheet.Range("D2").Value = "'1.87"
$newSheet.Range("E2").Value = "'91.11"
$newSheet.Range("F2").Value = "'1.47"
$newSheet.Range("G2").Value = "'0.0275"
$newSheet.Range("H2").Value = 7

# Data row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'009727"
$newSheet.Range("C3").Value = "招商中证500等权重指数增强C"
$newSheet.Range("D3").Value = "'0.69"
$newSheet.Range("E3").Value = "'91.11"
$newSheet.Range("F3").Value = "'1.47"
$newSheet.Range("G3").Value = "'0.0101"
$newSheet.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to "总计", pushing the existing
#    "2021-Q3" row down to row 3.
#    NOTE: re-resolve the "总计" worksheet by name — after Copy() above,
#    the original $totalSheet variable now refers to the newly-created
#    duplicate sheet, not the original "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Give the new index cell (A2) the same style as the header/index column.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.04

# Fix up the (now shifted) original row's index value.
$totalSheet.Range("A3").Value = 1
